# Strike through the six "Optimizing the game" bullet points that cover
# audio, mesh and texture performance work:
#   - Deactivate Read/write on textures and on models that don't need mesh colliders
#   - Disable rig on non-character models
#   - Enable Mesh compression (except vertex) (models)
#   - Mesh Renderer - Cast Shadow - off, Receive Shadow - off, use light probes - off, Reflection probes
#   - Ensure sizes aren't too large for textures, 1024 x 1024 UI atlases, 512 x 512 model textures
#   - AUDIO - Vorbis compression for Android, "Force Mono", Set Bitrate as low as possible

$d = $word.ActiveDocument

$firstText = "Deactivate Read/write on textures and on models that don't need mesh colliders"
$lastText  = "AUDIO - Vorbis compression for Android, ""Force Mono"", Set Bitrate as low as possible"

$firstIndex = 0
$lastIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $firstText) { $firstIndex = $i }
    if ($t -eq $lastText) { $lastIndex = $i }
}

# Fallback to the known layout of "To-Do List.docx" in case the text search
# above didn't find a unique match (e.g. whitespace differences).
if ($firstIndex -eq 0) { $firstIndex = 30 }
if ($lastIndex -eq 0) { $lastIndex = 35 }

# The boundary paragraphs only get their run text struck through (selection
# doesn't extend across their own paragraph mark), the paragraphs fully in
# between get both their run text and paragraph mark struck through.
$p = $d.Paragraphs.Item($firstIndex)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Font.StrikeThrough = $true

for ($i = $firstIndex + 1; $i -le $lastIndex - 1; $i++) {
    $d.Paragraphs.Item($i).Range.Font.StrikeThrough = $true
}

$p = $d.Paragraphs.Item($lastIndex)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Font.StrikeThrough = $true
